# Fruta / hortaliza, semanal
# Update Fecha (D), Calidad (I), Volumen (J), Precio minimo (K),
# Precio maximo (L), Precio promedio ponderado (M) and Precio $/Kg (P)
# for the weekly Repollo dataset (rows 2-21).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row: Fecha, Calidad, Volumen, Precio minimo, Precio maximo, Precio promedio ponderado, Precio $/Kg
$rows = @{
  2  = @(44474, "Segunda", 200,  600,  700,  650,  650)
  3  = @(44210, "Segunda", 900,  600,  700,  650,  650)
  4  = @(44544, "Primera", 1000, 600,  650,  625,  625)
  5  = @(44253, "Segunda", 1000, 800,  900,  850,  850)
  6  = @(44253, "Tercera", 800,  600,  700,  650,  650)
  7  = @(44174, "Segunda", 800,  450,  500,  475,  475)
  8  = @(44174, "Tercera", 1200, 250,  350,  300,  300)
  9  = @(44267, "Tercera", 400,  500,  600,  550,  550)
  10 = @(44573, "Tercera", 800,  600,  650,  625,  625)
  11 = @(44201, "Segunda", 500,  800,  900,  850,  850)
  12 = @(44874, "Tercera", 1200, 450,  500,  475,  475)
  13 = @(44278, "Segunda", 700,  600,  700,  650,  650)
  14 = @(44278, "Tercera", 400,  500,  600,  550,  550)
  15 = @(44245, "Primera", 800,  850,  900,  875,  875)
  16 = @(44245, "Segunda", 1000, 750,  800,  775,  775)
  17 = @(44799, "Primera", 800,  1000, 1200, 1100, 1100)
  18 = @(44935, "Segunda", 1000, 400,  500,  460,  460)
  19 = @(44658, "Segunda", 1000, 600,  650,  625,  625)
  20 = @(44224, "Segunda", 800,  850,  900,  875,  875)
  21 = @(44229, "Segunda", 760,  550,  600,  575,  575)
}

foreach ($r in $rows.Keys) {
  $vals = $rows[$r]
  $ws.Range("D$r").Value = $vals[0]
  $ws.Range("I$r").Value = $vals[1]
  $ws.Range("J$r").Value = $vals[2]
  $ws.Range("K$r").Value = $vals[3]
  $ws.Range("L$r").Value = $vals[4]
  $ws.Range("M$r").Value = $vals[5]
  $ws.Range("P$r").Value = $vals[6]
}
